$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 with new values
$ws.Range("B74").Value = 1206
$ws.Range("C74").Value = 1461
$ws.Range("D74").Value = -320
$ws.Range("G74").Value = 1414
$ws.Range("I74").Value = 2317
$ws.Range("J74").Value = 75
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 2242

# Add new row 75
# (Entered via a text formula + paste-values so Excel stores the
# date-like label "01-04-2021" as plain text, matching the other
# period-label cells in column A, instead of auto-converting it to a
# date serial number.)
$ws.Range("A75").Formula = '="01-04-2021"'
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("B75").Value = 1284
$ws.Range("C75").Value = 412
$ws.Range("D75").Value = -401
$ws.Range("E75").Value = -8347
$ws.Range("F75").Value = 14
$ws.Range("G75").Value = -7039
$ws.Range("H75").Value = -8347
$ws.Range("I75").Value = 1308
$ws.Range("J75").Value = 79
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 1229
